# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2404"
#   "<name>_new" -> "<name>_FV2410"
# and turn the sheet's used range into a real Excel Table, with the header
# row frozen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-J: headers that used to end in "_old" now end in "_FV2404"
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

# Columns L-U: headers that used to end in "_new" now end in "_FV2410"
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Headers[$i]
}

# Column K (11) is the "diff" column and keeps its name.

for ($i = 0; $i -lt $fv2410Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2410Headers[$i]
}

# Turn A1:U57 into a proper Excel Table named "Table1" (adds autoFilter +
# tableColumns matching the renamed headers).
$tableRange = $ws.Range("A1:U57")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
